# Updated symbol list on Fri Feb 10 08:54:50 UTC 2023 with GitHub Actions
# Applies refreshed Price (column D) and Volume(1h) (column E) values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row => @{ D = newPrice; E = newVolume }  (omit key if that column is unchanged)
$updates = @{
    2  = @{ D = "307.86";    E = "-4.62%" }
    3  = @{ D = "40.06";     E = "-6.74%" }
    4  = @{ D = "5.131";     E = "-0.99%" }
    5  = @{ D = "0.07744";   E = "-5.51%" }
    6  = @{ D = "4.248" }
    7  = @{ D = "1.621";     E = "-11.55%" }
    8  = @{ D = "0.8808";    E = "-5.52%" }
    9  = @{ D = "0.1011";    E = "-9.27%" }
    10 = @{ D = "0.1749";    E = "-6.46%" }
    11 = @{ D = "0.08964";   E = "-5.62%" }
    12 = @{ D = "0.04399";   E = "-4.88%" }
    13 = @{ E = "-0.22%" }
    14 = @{ D = "0.001258";  E = "-2.80%" }
    15 = @{ D = "0.005784";  E = "-1.15%" }
    16 = @{ E = "-0.27%" }
    18 = @{ D = "0.3322";    E = "-1.60%" }
    19 = @{ D = "7.002";     E = "-5.75%" }
    20 = @{ E = "-3.52%" }
    21 = @{ D = "0.2842";    E = "14.07%" }
    22 = @{ D = "0.04171";   E = "0.24%" }
    23 = @{ D = "0.001201";  E = "-3.58%" }
    24 = @{ D = "0.004105";  E = "-5.51%" }
    25 = @{ D = "0.0001300"; E = "8.40%" }
    26 = @{ E = "0.14%" }
    38 = @{ D = "0.02359";   E = "-14.36%" }
    39 = @{ D = "0.05144";   E = "-7.93%" }
    40 = @{ D = "0.007948";  E = "-4.48%" }
    41 = @{ E = "-5.00%" }
    42 = @{ D = "0.006387";  E = "-2.29%" }
    43 = @{ D = "0.001978";  E = "-5.38%" }
    44 = @{ D = "0.008889";  E = "18.33%" }
    45 = @{ D = "0.3332";    E = "-4.79%" }
    46 = @{ D = "0.00006564"; E = "-6.15%" }
    47 = @{ E = "0.06%" }
    48 = @{ E = "98.58%" }
    49 = @{ D = "0.006366";  E = "83.04%" }
    50 = @{ D = "0.00002100"; E = "0.06%" }
    51 = @{ D = "0.0002000"; E = "0.06%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    if ($cols.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["D"]
    }
    if ($cols.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cols["E"]
    }
}
